# "add page scope for watch number of registration"
# The registration sheet listed two registrants (rows 2-3) plus a stray
# timestamp value living on its own in row 4 (A4). The edit narrows the
# sheet down to a single registrant (row 2) and moves the timestamp up
# into row 3 / column A (updating it along the way), dropping every other
# field that used to be on row 3 and removing row 4 entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the second registrant's details (Суханов Игорь Константинович, ...)
# from row 3, keeping only column A for the timestamp that used to sit in A4.
$ws.Range("B3:I3").ClearContents()

# Row 4 (which only held the "2019-08-03 00:00:00" watch timestamp) is
# removed entirely, shrinking the used range to A1:I3.
$ws.Rows("4:4").Delete()

# Re-home the watch timestamp on row 3 with its updated value.
$ws.Range("A3").Value = "2019-09-18 00:00:00"
